$d = $word.ActiveDocument

# --- Paragraph 10: "Ernesto Canio ..." ---
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range

$xml10 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Ernesto Canio &#8211; Usuario AWS</w:t></w:r>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> &#34;ernesto0369&#34;</w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r10.InsertXML($xml10)

# --- Paragraph 11: "Usuario Arduino cloud ..." -> Ivan Cruz Canio ---
$p11 = $d.Paragraphs(11)
$r11 = $p11.Range

$xml11 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Ivan</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Cruz Canio &#8211; Usuario AWS</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:&#34;ivanc</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>2111&#34;</w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r11.InsertXML($xml11)

# --- Paragraph 12: remove the bookmark (now empty paragraph) ---
$p12 = $d.Paragraphs(12)
$r12 = $p12.Range
$xml12 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r12.InsertXML($xml12)
